$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.956999999999998
$ws.Range("C3").Value = -11.46959999999999
$ws.Range("A4").Value = -21.30400000000001
$ws.Range("B4").Value = 4.727000000000005
$ws.Range("C4").Value = -11.02739999999999
$ws.Range("B5").Value = 5.588199999999997
$ws.Range("D5").Value = -8.707899999999993
$ws.Range("A6").Value = -21.79050000000002
$ws.Range("B6").Value = 5.688799999999997
$ws.Range("A7").Value = -21.14910000000001
$ws.Range("A8").Value = -21.52830000000002
$ws.Range("B8").Value = 4.682700000000002
$ws.Range("C9").Value = -11.5019
$ws.Range("C11").Value = -13.7153
$ws.Range("C14").Value = -11.63539999999999
$ws.Range("A16").Value = -21.48610000000003
$ws.Range("B16").Value = 5.206799999999999
$ws.Range("C18").Value = -14.52500000000002
$ws.Range("A20").Value = -22.85380000000002
$ws.Range("D20").Value = -8.381400000000006
$ws.Range("A21").Value = -20.34909999999999
$ws.Range("B22").Value = 6.123299999999993
$ws.Range("C25").Value = -11.36399999999999
